# Correcciones en general y todotid
# Add two new rows to the "edit-fields" sheet describing the new
# "cart" node's "session_id" field (its relation + its type).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("edit-fields")

$ws.Cells.Item(84, 1).Value = "cart"
$ws.Cells.Item(84, 2).Value = "session_id"
$ws.Cells.Item(84, 3).Value = "relation"
$ws.Cells.Item(84, 4).Value = 0

$ws.Cells.Item(85, 1).Value = "cart"
$ws.Cells.Item(85, 2).Value = "session_id"
$ws.Cells.Item(85, 3).Value = "type"
$ws.Cells.Item(85, 4).Value = "string"

# Match the author's final selection/active cell on this sheet.
$ws.Range("B85").Select() | Out-Null
